$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$dcell = $ws.Range("D2")
$dcell.NumberFormat = "@"
$dcell.Value = '30.086.57'
$dcell.Style = "Normal"
$ws.Range("E2").Value = '  -1.41%  '

$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$dcell = $ws.Range("D3")
$dcell.NumberFormat = "@"
$dcell.Value = '1.850.38'
$dcell.Style = "Normal"
$ws.Range("E3").Value = '  -0.14%  '

$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$dcell = $ws.Range("D4")
$dcell.NumberFormat = "@"
$dcell.Value = '1.002'
$dcell.Style = "Normal"
$ws.Range("E4").Value = '  +0.14%  '

$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$dcell = $ws.Range("D5")
$dcell.NumberFormat = "@"
$dcell.Value = '235.77'
$dcell.Style = "Normal"
$ws.Range("E5").Value = '  +0.83%  '

$ws.Range("B6").Value = 'USDC'
$ws.Range("C6").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$dcell = $ws.Range("D6")
$dcell.NumberFormat = "@"
$dcell.Value = '1.002'
$dcell.Style = "Normal"
$ws.Range("E6").Value = '  +0.14%  '

$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$dcell = $ws.Range("D7")
$dcell.NumberFormat = "@"
$dcell.Value = '0.4639'
$dcell.Style = "Normal"
$ws.Range("E7").Value = '  -1.43%  '

$ws.Range("B8").Value = 'OKB'
$ws.Range("C8").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$dcell = $ws.Range("D8")
$dcell.NumberFormat = "@"
$dcell.Value = '43.98'
$dcell.Style = "Normal"
$ws.Range("E8").Value = '  +1.69%  '

$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$dcell = $ws.Range("D9")
$dcell.NumberFormat = "@"
$dcell.Value = '0.2773'
$dcell.Style = "Normal"
$ws.Range("E9").Value = '  +1.08%  '

$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$dcell = $ws.Range("D10")
$dcell.NumberFormat = "@"
$dcell.Value = '0.06389'
$dcell.Style = "Normal"
$ws.Range("E10").Value = '  +0.79%  '

$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$dcell = $ws.Range("D11")
$dcell.NumberFormat = "@"
$dcell.Value = '18.11'
$dcell.Style = "Normal"
$ws.Range("E11").Value = '  +3.23%  '

$ws.Range("B12").Value = 'Litecoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$dcell = $ws.Range("D12")
$dcell.NumberFormat = "@"
$dcell.Value = '96.52'
$dcell.Style = "Normal"
$ws.Range("E12").Value = '  +13.98%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$dcell = $ws.Range("D13")
$dcell.NumberFormat = "@"
$dcell.Value = '1.842.20'
$dcell.Style = "Normal"
$ws.Range("E13").Value = '  +0.51%  '

$ws.Range("B14").Value = 'TRON'
$ws.Range("C14").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$dcell = $ws.Range("D14")
$dcell.NumberFormat = "@"
$dcell.Value = '0.07524'
$dcell.Style = "Normal"
$ws.Range("E14").Value = '  +1.32%  '

$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$dcell = $ws.Range("D15")
$dcell.NumberFormat = "@"
$dcell.Value = '4.967'
$dcell.Style = "Normal"
$ws.Range("E15").Value = '  -1.78%  '

$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$dcell = $ws.Range("D16")
$dcell.NumberFormat = "@"
$dcell.Value = '0.6241'
$dcell.Style = "Normal"
$ws.Range("E16").Value = '  -0.54%  '

$ws.Range("B17").Value = 'BitcoinCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$dcell = $ws.Range("D17")
$dcell.NumberFormat = "@"
$dcell.Value = '293.25'
$dcell.Style = "Normal"
$ws.Range("E17").Value = '  +20.04%  '

$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$dcell = $ws.Range("D18")
$dcell.NumberFormat = "@"
$dcell.Value = '30.085.51'
$dcell.Style = "Normal"
$ws.Range("E18").Value = '  -1.36%  '

$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$dcell = $ws.Range("D19")
$dcell.NumberFormat = "@"
$dcell.Value = '1.001'
$dcell.Style = "Normal"
$ws.Range("E19").Value = '  +0.04%  '

$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$dcell = $ws.Range("D20")
$dcell.NumberFormat = "@"
$dcell.Value = '12.66'
$dcell.Style = "Normal"
$ws.Range("E20").Value = '  -0.34%  '

$ws.Range("B21").Value = 'ShibaInu'
$ws.Range("C21").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$dcell = $ws.Range("D21")
$dcell.NumberFormat = "@"
$dcell.Value = '0.000007362'
$dcell.Style = "Normal"
$ws.Range("E21").Value = '  +0.26%  '

$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$dcell = $ws.Range("D22")
$dcell.NumberFormat = "@"
$dcell.Value = '2.085.17'
$dcell.Style = "Normal"
$ws.Range("E22").Value = '  -0.28%  '

$ws.Range("B23").Value = 'BinanceUSD'
$ws.Range("C23").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$dcell = $ws.Range("D23")
$dcell.NumberFormat = "@"
$dcell.Value = '1.003'
$dcell.Style = "Normal"
$ws.Range("E23").Value = '  +0.25%  '

$ws.Range("B24").Value = 'Uniswap'
$ws.Range("C24").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$dcell = $ws.Range("D24")
$dcell.NumberFormat = "@"
$dcell.Value = '4.981'
$dcell.Style = "Normal"
$ws.Range("E24").Value = '  +0.36%  '

$ws.Range("B25").Value = 'Chainlink'
$ws.Range("C25").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$dcell = $ws.Range("D25")
$dcell.NumberFormat = "@"
$dcell.Value = '6.099'
$dcell.Style = "Normal"
$ws.Range("E25").Value = '  +1.63%  '

$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$dcell = $ws.Range("D26")
$dcell.NumberFormat = "@"
$dcell.Value = '165.01'
$dcell.Style = "Normal"
$ws.Range("E26").Value = '  +1.57%  '

$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$dcell = $ws.Range("D27")
$dcell.NumberFormat = "@"
$dcell.Value = '9.065'
$dcell.Style = "Normal"
$ws.Range("E27").Value = '  -2.31%  '

$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$dcell = $ws.Range("D28")
$dcell.NumberFormat = "@"
$dcell.Value = '19.20'
$dcell.Style = "Normal"
$ws.Range("E28").Value = '  +6.20%  '

$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$dcell = $ws.Range("D29")
$dcell.NumberFormat = "@"
$dcell.Value = '1.934'
$dcell.Style = "Normal"
$ws.Range("E29").Value = '  +2.64%  '

$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$dcell = $ws.Range("D30")
$dcell.NumberFormat = "@"
$dcell.Value = '0.1072'
$dcell.Style = "Normal"
$ws.Range("E30").Value = '  +6.20%  '

$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$dcell = $ws.Range("D31")
$dcell.NumberFormat = "@"
$dcell.Value = '1.325'
$dcell.Style = "Normal"
$ws.Range("E31").Value = '  -3.30%  '

$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$dcell = $ws.Range("D32")
$dcell.NumberFormat = "@"
$dcell.Value = '4.000'
$dcell.Style = "Normal"
$ws.Range("E32").Value = '  -0.93%  '

$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$dcell = $ws.Range("D33")
$dcell.NumberFormat = "@"
$dcell.Value = '3.811'
$dcell.Style = "Normal"
$ws.Range("E33").Value = '  -1.17%  '

$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$dcell = $ws.Range("D34")
$dcell.NumberFormat = "@"
$dcell.Value = '0.04899'
$dcell.Style = "Normal"
$ws.Range("E34").Value = '  -0.05%  '

$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$dcell = $ws.Range("D35")
$dcell.NumberFormat = "@"
$dcell.Value = '0.7266'
$dcell.Style = "Normal"
$ws.Range("E35").Value = '  +3.03%  '

$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$dcell = $ws.Range("D36")
$dcell.NumberFormat = "@"
$dcell.Value = '1.109'
$dcell.Style = "Normal"
$ws.Range("E36").Value = '  -2.47%  '

$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$dcell = $ws.Range("D37")
$dcell.NumberFormat = "@"
$dcell.Value = '2.727'
$dcell.Style = "Normal"
$ws.Range("E37").Value = '  +0.81%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$dcell = $ws.Range("D38")
$dcell.NumberFormat = "@"
$dcell.Value = '0.01892'
$dcell.Style = "Normal"
$ws.Range("E38").Value = '  -0.91%  '

$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$dcell = $ws.Range("D39")
$dcell.NumberFormat = "@"
$dcell.Value = '2.656'
$dcell.Style = "Normal"
$ws.Range("E39").Value = '  -1.08%  '

$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$dcell = $ws.Range("D40")
$dcell.NumberFormat = "@"
$dcell.Value = '1.964'
$dcell.Style = "Normal"
$ws.Range("E40").Value = '  -0.64%  '

$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$dcell = $ws.Range("D41")
$dcell.NumberFormat = "@"
$dcell.Value = '0.8581'
$dcell.Style = "Normal"
$ws.Range("E41").Value = '  -1.87%  '

$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$dcell = $ws.Range("D42")
$dcell.NumberFormat = "@"
$dcell.Value = '104.58'
$dcell.Style = "Normal"
$ws.Range("E42").Value = '  -0.60%  '

$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$dcell = $ws.Range("D43")
$dcell.NumberFormat = "@"
$dcell.Value = '1.001'
$dcell.Style = "Normal"
$ws.Range("E43").Value = '  +0.09%  '

$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$dcell = $ws.Range("D44")
$dcell.NumberFormat = "@"
$dcell.Value = '5.677'
$dcell.Style = "Normal"
$ws.Range("E44").Value = '  +3.18%  '

$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$dcell = $ws.Range("D45")
$dcell.NumberFormat = "@"
$dcell.Value = '0.4031'
$dcell.Style = "Normal"
$ws.Range("E45").Value = '  -1.13%  '

$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$dcell = $ws.Range("D46")
$dcell.NumberFormat = "@"
$dcell.Value = '65.09'
$dcell.Style = "Normal"
$ws.Range("E46").Value = '  +3.38%  '

$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$dcell = $ws.Range("D47")
$dcell.NumberFormat = "@"
$dcell.Value = '7.031'
$dcell.Style = "Normal"
$ws.Range("E47").Value = '  -3.06%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$dcell = $ws.Range("D48")
$dcell.NumberFormat = "@"
$dcell.Value = '8.957'
$dcell.Style = "Normal"
$ws.Range("E48").Value = '  +4.45%  '

$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$dcell = $ws.Range("D49")
$dcell.NumberFormat = "@"
$dcell.Value = '0.1185'
$dcell.Style = "Normal"
$ws.Range("E49").Value = '  -1.17%  '

$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$dcell = $ws.Range("D50")
$dcell.NumberFormat = "@"
$dcell.Value = '33.85'
$dcell.Style = "Normal"
$ws.Range("E50").Value = '  +1.34%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$dcell = $ws.Range("D51")
$dcell.NumberFormat = "@"
$dcell.Value = '0.05543'
$dcell.Style = "Normal"
$ws.Range("E51").Value = '  +0.21%  '
